$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Split the "Application Description" game paragraph: append two new
#    sentences to the end of the first paragraph (about Tableau piles), and
#    expand the second "Game features" paragraph with extra clauses.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Tableau area consisting of thirteen card piles.") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" The game will only allow the user to make moves that are within the parameters of the game. The game will also have simple and clear graphical elements.")

$rng = $d.Content
$rng.Find.Execute("Other features such as a move counter or a timer may be added to the game as project scheduling allows.") | Out-Null
$rng.Text = "Additional features such as a move counter, a timer, an undo move button, or a restart current game button may be added to the game as project scheduling allows. Other advanced features such as multiplayer or online leaderboards are outside of the scope of this project."

# ---------------------------------------------------------------------------
# 2. "Deliverables in this project will submitted in four main areas:"
#    -> "...will be submitted..." (insert the missing "be ")
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Deliverables in this project will ") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("be ")

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the document to right after
#    the colon that ends the "...four main areas:" paragraph. We do this by
#    temporarily inserting a placeholder character after the colon, wrapping
#    a bookmark around it (non-collapsed ranges anchor correctly), and then
#    deleting the placeholder while the (now zero-width) bookmark remains in
#    place.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("submitted in four main areas:") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("X")
$phRng = $d.Range($rng.Start, $rng.Start + 1)
$d.Bookmarks.Add("_GoBack", $phRng)
$phRng = $d.Range($rng.Start, $rng.Start + 1)
$phRng.Text = ""

# ---------------------------------------------------------------------------
# 4. Add trailing periods to several bullet items that were missing them.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Includes high-level design information, conceptual system design information, and technical design information") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

$rng = $d.Content
$rng.Find.Execute("Additional Documents – Consists of a Test Specification document, team time sheets, and meeting minutes") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

$rng = $d.Content
$rng.Find.Execute("Ensure that all data is backed up to an off-site location") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

$rng = $d.Content
$rng.Find.Execute("Ensure that all data is accessible by the team at any time") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

$rng = $d.Content
$rng.Find.Execute("Allow for all changes to documents to be recorded") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

$rng = $d.Content
$rng.Find.Execute("Allow any changes to be rolled back") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

# ---------------------------------------------------------------------------
# 5. Insert a new sentence before the GitHub's archiving/recordkeeping
#    sentence at the end of the Data Management Plan section.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("GitHub’s archiving and recordkeeping tools mentioned above") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("All documents are readily available for revision for each team member. ")
